# Generate Report for Handoff
# Update the handoff/handback timestamps for e602f38e-9f32-4884-acc3-9b89f3f87e67.md
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row 5 is the e602f38e-9f32-4884-acc3-9b89f3f87e67.md file,
# column D = "Latest Handoff Date"
$wsOverview.Range("D5").Value = "2016-03-31 05:02:23"

# zh-cn sheet: row 5 is the e602f38e-9f32-4884-acc3-9b89f3f87e67.md file,
# column E = "Latest Handoff Datetime"
$wsZhCn.Range("E5").Value = "2016-03-31 05:02:13"

# de-de sheet: row 5 is the e602f38e-9f32-4884-acc3-9b89f3f87e67.md file,
# column E = "Latest Handoff Datetime"
$wsDeDe.Range("E5").Value = "2016-03-31 05:02:23"
